$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at position 595, shifting existing rows 595:683 down to 597:685
$ws.Range("A595:R596").EntireRow.Insert()

# Fill new row 595 with the latest weekly price record
$ws.Cells.Item(595, 1).Value = 9
$ws.Cells.Item(595, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(595, 3).Value = "Metropolitana"
$ws.Cells.Item(595, 4).Value = 44776
$ws.Cells.Item(595, 5).Value = 13
$ws.Cells.Item(595, 6).Value = 100112040
$ws.Cells.Item(595, 7).Value = "Cilantro"
$ws.Cells.Item(595, 8).Value = "Sin especificar"
$ws.Cells.Item(595, 9).Value = "Primera"
$ws.Cells.Item(595, 10).Value = 70
$ws.Cells.Item(595, 11).Value = 9000
$ws.Cells.Item(595, 12).Value = 9000
$ws.Cells.Item(595, 13).Value = 9000
$ws.Cells.Item(595, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(595, 15).Value = "Región Metropolitana"
$ws.Cells.Item(595, 16).Value = 250
$ws.Cells.Item(595, 17).Value = 36
$ws.Cells.Item(595, 18).Value = "Hortaliza"

# Fill new row 596 with the latest weekly price record (second unit type)
$ws.Cells.Item(596, 1).Value = 9
$ws.Cells.Item(596, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(596, 3).Value = "Metropolitana"
$ws.Cells.Item(596, 4).Value = 44776
$ws.Cells.Item(596, 5).Value = 13
$ws.Cells.Item(596, 6).Value = 100112040
$ws.Cells.Item(596, 7).Value = "Cilantro"
$ws.Cells.Item(596, 8).Value = "Sin especificar"
$ws.Cells.Item(596, 9).Value = "Primera"
$ws.Cells.Item(596, 10).Value = 160
$ws.Cells.Item(596, 11).Value = 15000
$ws.Cells.Item(596, 12).Value = 16000
$ws.Cells.Item(596, 13).Value = 15500
$ws.Cells.Item(596, 14).Value = "$/docena de atados"
$ws.Cells.Item(596, 15).Value = "Región Metropolitana"
$ws.Cells.Item(596, 16).Value = 5167
$ws.Cells.Item(596, 17).Value = 3
$ws.Cells.Item(596, 18).Value = "Hortaliza"
